# Daten aktualisiert am 2023-09-24
# Refresh the cryptocurrency market snapshot on Sheet1 (rows 2-51: Ticker, Name,
# Price, Market Cap, Volume, Change (24h)). Some coins changed rank order because
# the sheet is sorted by Market Cap (column E) descending, so a handful of rows'
# Ticker/Name also had to move to stay next to their refreshed numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 26604  # D2 Price
$ws.Cells.Item(2, 5).Value = 518561249912  # E2 Market Cap
$ws.Cells.Item(2, 6).Value = 5286209623  # F2 Volume
$ws.Cells.Item(2, 7).Value = 0.12943  # G2 Change (24h)

$ws.Cells.Item(3, 4).Value = 1593.69  # D3 Price
$ws.Cells.Item(3, 5).Value = 191627710580  # E3 Market Cap
$ws.Cells.Item(3, 6).Value = 2651407000  # F3 Volume
$ws.Cells.Item(3, 7).Value = 0.07874  # G3 Change (24h)

$ws.Cells.Item(4, 4).Value = 0.999966  # D4 Price
$ws.Cells.Item(4, 5).Value = 83203330822  # E4 Market Cap
$ws.Cells.Item(4, 6).Value = 8061869039  # F4 Volume
$ws.Cells.Item(4, 7).Value = 0.02274  # G4 Change (24h)

$ws.Cells.Item(5, 4).Value = 210.8  # D5 Price
$ws.Cells.Item(5, 5).Value = 32437149245  # E5 Market Cap
$ws.Cells.Item(5, 6).Value = 173764446  # F5 Volume
$ws.Cells.Item(5, 7).Value = 0.03377  # G5 Change (24h)

$ws.Cells.Item(6, 4).Value = 0.510581  # D6 Price
$ws.Cells.Item(6, 5).Value = 27212561438  # E6 Market Cap
$ws.Cells.Item(6, 6).Value = 313455140  # F6 Volume
$ws.Cells.Item(6, 7).Value = -0.16132  # G6 Change (24h)

$ws.Cells.Item(7, 3).Value = "USDC"  # C7 Name
$ws.Cells.Item(7, 4).Value = 1  # D7 Price
$ws.Cells.Item(7, 5).Value = 25774759109  # E7 Market Cap
$ws.Cells.Item(7, 6).Value = 1600202498  # F7 Volume
$ws.Cells.Item(7, 7).Value = 0.02076  # G7 Change (24h)

$ws.Cells.Item(8, 4).Value = 1594.35  # D8 Price
$ws.Cells.Item(8, 5).Value = 13911269449  # E8 Market Cap
$ws.Cells.Item(8, 6).Value = 7147333  # F8 Volume
$ws.Cells.Item(8, 7).Value = 0.17159  # G8 Change (24h)

$ws.Cells.Item(9, 4).Value = 0.061518  # D9 Price
$ws.Cells.Item(9, 5).Value = 8685449215  # E9 Market Cap
$ws.Cells.Item(9, 6).Value = 109822523  # F9 Volume
$ws.Cells.Item(9, 7).Value = -0.13686  # G9 Change (24h)

$ws.Cells.Item(10, 4).Value = 0.245805  # D10 Price
$ws.Cells.Item(10, 5).Value = 8616299607  # E10 Market Cap
$ws.Cells.Item(10, 6).Value = 60362852  # F10 Volume
$ws.Cells.Item(10, 7).Value = -0.22022  # G10 Change (24h)

$ws.Cells.Item(11, 2).Value = "SOL"  # B11 Ticker
$ws.Cells.Item(11, 3).Value = "Solana"  # C11 Name
$ws.Cells.Item(11, 4).Value = 19.62  # D11 Price
$ws.Cells.Item(11, 5).Value = 8106142524  # E11 Market Cap
$ws.Cells.Item(11, 6).Value = 111207653  # F11 Volume
$ws.Cells.Item(11, 7).Value = 0.75925  # G11 Change (24h)

$ws.Cells.Item(12, 2).Value = "TON"  # B12 Ticker
$ws.Cells.Item(12, 3).Value = "Toncoin"  # C12 Name
$ws.Cells.Item(12, 4).Value = 2.23  # D12 Price
$ws.Cells.Item(12, 5).Value = 7703772619  # E12 Market Cap
$ws.Cells.Item(12, 6).Value = 12353611  # F12 Volume
$ws.Cells.Item(12, 7).Value = -4.14935  # G12 Change (24h)

$ws.Cells.Item(13, 4).Value = 0.084181  # D13 Price
$ws.Cells.Item(13, 5).Value = 7506255995  # E13 Market Cap
$ws.Cells.Item(13, 6).Value = 135022365  # F13 Volume
$ws.Cells.Item(13, 7).Value = 0.77325  # G13 Change (24h)

$ws.Cells.Item(14, 4).Value = 4.04  # D14 Price
$ws.Cells.Item(14, 5).Value = 5168460808  # E14 Market Cap
$ws.Cells.Item(14, 6).Value = 54035810  # F14 Volume
$ws.Cells.Item(14, 7).Value = 0.84137  # G14 Change (24h)

$ws.Cells.Item(15, 4).Value = 0.52183  # D15 Price
$ws.Cells.Item(15, 5).Value = 4862313303  # E15 Market Cap
$ws.Cells.Item(15, 6).Value = 73001878  # F15 Volume
$ws.Cells.Item(15, 7).Value = 0.00797  # G15 Change (24h)

$ws.Cells.Item(16, 4).Value = 64.66  # D16 Price
$ws.Cells.Item(16, 5).Value = 4766430982  # E16 Market Cap
$ws.Cells.Item(16, 6).Value = 184268355  # F16 Volume
$ws.Cells.Item(16, 7).Value = -0.23203  # G16 Change (24h)

$ws.Cells.Item(17, 2).Value = "SHIB"  # B17 Ticker
$ws.Cells.Item(17, 3).Value = "Shiba Inu"  # C17 Name
$ws.Cells.Item(17, 4).Value = 0.00000737  # D17 Price
$ws.Cells.Item(17, 5).Value = 4339994115  # E17 Market Cap
$ws.Cells.Item(17, 6).Value = 151783928  # F17 Volume
$ws.Cells.Item(17, 7).Value = 0.18089  # G17 Change (24h)

$ws.Cells.Item(18, 2).Value = "WBTC"  # B18 Ticker
$ws.Cells.Item(18, 3).Value = "Wrapped Bitcoin"  # C18 Name
$ws.Cells.Item(18, 4).Value = 26567  # D18 Price
$ws.Cells.Item(18, 5).Value = 4326498996  # E18 Market Cap
$ws.Cells.Item(18, 6).Value = 18358826  # F18 Volume
$ws.Cells.Item(18, 7).Value = 0.12348  # G18 Change (24h)

$ws.Cells.Item(19, 4).Value = 208.64  # D19 Price
$ws.Cells.Item(19, 5).Value = 4073114887  # E19 Market Cap
$ws.Cells.Item(19, 6).Value = 71390949  # F19 Volume
$ws.Cells.Item(19, 7).Value = 0.30936  # G19 Change (24h)

$ws.Cells.Item(20, 2).Value = "LINK"  # B20 Ticker
$ws.Cells.Item(20, 3).Value = "Chainlink"  # C20 Name
$ws.Cells.Item(20, 4).Value = 7.1  # D20 Price
$ws.Cells.Item(20, 5).Value = 3951101258  # E20 Market Cap
$ws.Cells.Item(20, 6).Value = 259440126  # F20 Volume
$ws.Cells.Item(20, 7).Value = 1.30578  # G20 Change (24h)

$ws.Cells.Item(21, 2).Value = "DAI"  # B21 Ticker
$ws.Cells.Item(21, 3).Value = "Dai"  # C21 Name
$ws.Cells.Item(21, 4).Value = 1  # D21 Price
$ws.Cells.Item(21, 5).Value = 3843467995  # E21 Market Cap
$ws.Cells.Item(21, 6).Value = 57549519  # F21 Volume
$ws.Cells.Item(21, 7).Value = 0.09417  # G21 Change (24h)

$ws.Cells.Item(22, 2).Value = "TUSD"  # B22 Ticker
$ws.Cells.Item(22, 3).Value = "TrueUSD"  # C22 Name
$ws.Cells.Item(22, 4).Value = 0.99898  # D22 Price
$ws.Cells.Item(22, 5).Value = 3507279323  # E22 Market Cap
$ws.Cells.Item(22, 6).Value = 46760150  # F22 Volume
$ws.Cells.Item(22, 7).Value = 0.01398  # G22 Change (24h)

$ws.Cells.Item(23, 2).Value = "LEO"  # B23 Ticker
$ws.Cells.Item(23, 3).Value = "LEO Token"  # C23 Name
$ws.Cells.Item(23, 4).Value = 3.75  # D23 Price
$ws.Cells.Item(23, 5).Value = 3485171333  # E23 Market Cap
$ws.Cells.Item(23, 6).Value = 180450  # F23 Volume
$ws.Cells.Item(23, 7).Value = -2.26696  # G23 Change (24h)

$ws.Cells.Item(24, 2).Value = "UNI"  # B24 Ticker
$ws.Cells.Item(24, 3).Value = "Uniswap"  # C24 Name
$ws.Cells.Item(24, 4).Value = 4.28  # D24 Price
$ws.Cells.Item(24, 5).Value = 3224902193  # E24 Market Cap
$ws.Cells.Item(24, 6).Value = 39446281  # F24 Volume
$ws.Cells.Item(24, 7).Value = 0.1842  # G24 Change (24h)

$ws.Cells.Item(25, 2).Value = "AVAX"  # B25 Ticker
$ws.Cells.Item(25, 3).Value = "Avalanche"  # C25 Name
$ws.Cells.Item(25, 4).Value = 9.01  # D25 Price
$ws.Cells.Item(25, 5).Value = 3188662677  # E25 Market Cap
$ws.Cells.Item(25, 6).Value = 73583184  # F25 Volume
$ws.Cells.Item(25, 7).Value = 0.56853  # G25 Change (24h)

$ws.Cells.Item(26, 2).Value = "XLM"  # B26 Ticker
$ws.Cells.Item(26, 3).Value = "Stellar"  # C26 Name
$ws.Cells.Item(26, 4).Value = 0.113434  # D26 Price
$ws.Cells.Item(26, 5).Value = 3142869270  # E26 Market Cap
$ws.Cells.Item(26, 6).Value = 29612413  # F26 Volume
$ws.Cells.Item(26, 7).Value = -0.36688  # G26 Change (24h)

$ws.Cells.Item(27, 4).Value = 142.94  # D27 Price
$ws.Cells.Item(27, 5).Value = 2593165477  # E27 Market Cap
$ws.Cells.Item(27, 6).Value = 37268065  # F27 Volume
$ws.Cells.Item(27, 7).Value = -0.70748  # G27 Change (24h)

$ws.Cells.Item(28, 4).Value = 42.85  # D28 Price
$ws.Cells.Item(28, 5).Value = 2573393240  # E28 Market Cap
$ws.Cells.Item(28, 6).Value = 1437314  # F28 Volume
$ws.Cells.Item(28, 7).Value = 0.28249  # G28 Change (24h)

$ws.Cells.Item(29, 5).Value = 2381668298  # E29 Market Cap
$ws.Cells.Item(29, 6).Value = 684960026  # F29 Volume
$ws.Cells.Item(29, 7).Value = 0.00317  # G29 Change (24h)

$ws.Cells.Item(30, 4).Value = 15.27  # D30 Price
$ws.Cells.Item(30, 5).Value = 2184472345  # E30 Market Cap
$ws.Cells.Item(30, 6).Value = 38858988  # F30 Volume
$ws.Cells.Item(30, 7).Value = 0.47672  # G30 Change (24h)

$ws.Cells.Item(31, 4).Value = 7.06  # D31 Price
$ws.Cells.Item(31, 5).Value = 2064306033  # E31 Market Cap
$ws.Cells.Item(31, 6).Value = 79706843  # F31 Volume
$ws.Cells.Item(31, 7).Value = -0.44987  # G31 Change (24h)

$ws.Cells.Item(32, 4).Value = 0.050437  # D32 Price
$ws.Cells.Item(32, 5).Value = 1685403380  # E32 Market Cap
$ws.Cells.Item(32, 6).Value = 15665416  # F32 Volume
$ws.Cells.Item(32, 7).Value = -1.63668  # G32 Change (24h)

$ws.Cells.Item(33, 4).Value = 3.25  # D33 Price
$ws.Cells.Item(33, 5).Value = 1468511796  # E33 Market Cap
$ws.Cells.Item(33, 6).Value = 46443126  # F33 Volume
$ws.Cells.Item(33, 7).Value = 0.78301  # G33 Change (24h)

$ws.Cells.Item(34, 4).Value = 0.051313  # D34 Price
$ws.Cells.Item(34, 5).Value = 1353733630  # E34 Market Cap
$ws.Cells.Item(34, 6).Value = 13112024  # F34 Volume
$ws.Cells.Item(34, 7).Value = 1.58151  # G34 Change (24h)

$ws.Cells.Item(35, 4).Value = 1.49  # D35 Price
$ws.Cells.Item(35, 5).Value = 1323768554  # E35 Market Cap
$ws.Cells.Item(35, 6).Value = 9894668  # F35 Volume
$ws.Cells.Item(35, 7).Value = 0.26614  # G35 Change (24h)

$ws.Cells.Item(36, 4).Value = 2.96  # D36 Price
$ws.Cells.Item(36, 5).Value = 1315129194  # E36 Market Cap
$ws.Cells.Item(36, 6).Value = 9140714  # F36 Volume
$ws.Cells.Item(36, 7).Value = 0.35466  # G36 Change (24h)

$ws.Cells.Item(37, 4).Value = 90.11  # D37 Price
$ws.Cells.Item(37, 5).Value = 1311742154  # E37 Market Cap
$ws.Cells.Item(37, 6).Value = 9116989  # F37 Volume
$ws.Cells.Item(37, 7).Value = -0.70339  # G37 Change (24h)

$ws.Cells.Item(38, 4).Value = 0.39705  # D38 Price
$ws.Cells.Item(38, 5).Value = 1284065364  # E38 Market Cap
$ws.Cells.Item(38, 6).Value = 11973910  # F38 Volume
$ws.Cells.Item(38, 7).Value = 1.43755  # G38 Change (24h)

$ws.Cells.Item(39, 2).Value = "APT"  # B39 Ticker
$ws.Cells.Item(39, 3).Value = "Aptos"  # C39 Name
$ws.Cells.Item(39, 4).Value = 5.32  # D39 Price
$ws.Cells.Item(39, 5).Value = 1258831819  # E39 Market Cap
$ws.Cells.Item(39, 6).Value = 32581917  # F39 Volume
$ws.Cells.Item(39, 7).Value = 3.91372  # G39 Change (24h)

$ws.Cells.Item(40, 2).Value = "VET"  # B40 Ticker
$ws.Cells.Item(40, 3).Value = "VeChain"  # C40 Name
$ws.Cells.Item(40, 4).Value = 0.01691804  # D40 Price
$ws.Cells.Item(40, 5).Value = 1230279023  # E40 Market Cap
$ws.Cells.Item(40, 6).Value = 22117596  # F40 Volume
$ws.Cells.Item(40, 7).Value = 0.09121  # G40 Change (24h)

$ws.Cells.Item(41, 4).Value = 1284.43  # D41 Price
$ws.Cells.Item(41, 5).Value = 1157407724  # E41 Market Cap
$ws.Cells.Item(41, 6).Value = 34979828  # F41 Volume
$ws.Cells.Item(41, 7).Value = 0.49575  # G41 Change (24h)

$ws.Cells.Item(42, 2).Value = "ARB"  # B42 Ticker
$ws.Cells.Item(42, 3).Value = "Arbitrum"  # C42 Name
$ws.Cells.Item(42, 4).Value = 0.829186  # D42 Price
$ws.Cells.Item(42, 5).Value = 1056989634  # E42 Market Cap
$ws.Cells.Item(42, 6).Value = 40649559  # F42 Volume
$ws.Cells.Item(42, 7).Value = -0.04506  # G42 Change (24h)

$ws.Cells.Item(43, 4).Value = 1.11  # D43 Price
$ws.Cells.Item(43, 5).Value = 1046629964  # E43 Market Cap
$ws.Cells.Item(43, 6).Value = 24376532  # F43 Volume
$ws.Cells.Item(43, 7).Value = 0.02613  # G43 Change (24h)

$ws.Cells.Item(44, 2).Value = "OP"  # B44 Ticker
$ws.Cells.Item(44, 3).Value = "Optimism"  # C44 Name
$ws.Cells.Item(44, 4).Value = 1.29  # D44 Price
$ws.Cells.Item(44, 5).Value = 1033573587  # E44 Market Cap
$ws.Cells.Item(44, 6).Value = 32730120  # F44 Volume
$ws.Cells.Item(44, 7).Value = 0.29553  # G44 Change (24h)

$ws.Cells.Item(45, 4).Value = 0.04600801  # D45 Price
$ws.Cells.Item(45, 5).Value = 958455531  # E45 Market Cap
$ws.Cells.Item(45, 6).Value = 7400215  # F45 Volume
$ws.Cells.Item(45, 7).Value = 0.68474  # G45 Change (24h)

$ws.Cells.Item(46, 4).Value = 1734.01  # D46 Price
$ws.Cells.Item(46, 5).Value = 917043628  # E46 Market Cap
$ws.Cells.Item(46, 6).Value = 3699007  # F46 Volume
$ws.Cells.Item(46, 7).Value = 0.21567  # G46 Change (24h)

$ws.Cells.Item(47, 4).Value = 62.72  # D47 Price
$ws.Cells.Item(47, 5).Value = 914090373  # E47 Market Cap
$ws.Cells.Item(47, 6).Value = 38349036  # F47 Volume
$ws.Cells.Item(47, 7).Value = -0.8118  # G47 Change (24h)

$ws.Cells.Item(48, 4).Value = 0.087417  # D48 Price
$ws.Cells.Item(48, 5).Value = 807561413  # E48 Market Cap
$ws.Cells.Item(48, 6).Value = 14652153  # F48 Volume
$ws.Cells.Item(48, 7).Value = -1.38595  # G48 Change (24h)

$ws.Cells.Item(49, 2).Value = "ALGO"  # B49 Ticker
$ws.Cells.Item(49, 3).Value = "Algorand"  # C49 Name
$ws.Cells.Item(49, 4).Value = 0.101447  # D49 Price
$ws.Cells.Item(49, 5).Value = 794019115  # E49 Market Cap
$ws.Cells.Item(49, 6).Value = 15826954  # F49 Volume
$ws.Cells.Item(49, 7).Value = 0.33846  # G49 Change (24h)

$ws.Cells.Item(50, 2).Value = "WBT"  # B50 Ticker
$ws.Cells.Item(50, 3).Value = "WhiteBIT Coin"  # C50 Name
$ws.Cells.Item(50, 4).Value = 5.31  # D50 Price
$ws.Cells.Item(50, 5).Value = 764812777  # E50 Market Cap
$ws.Cells.Item(50, 6).Value = 7890381  # F50 Volume
$ws.Cells.Item(50, 7).Value = -0.14292  # G50 Change (24h)

$ws.Cells.Item(51, 2).Value = "USDD"  # B51 Ticker
$ws.Cells.Item(51, 3).Value = "USDD"  # C51 Name
$ws.Cells.Item(51, 4).Value = 1  # D51 Price
$ws.Cells.Item(51, 5).Value = 726432006  # E51 Market Cap
$ws.Cells.Item(51, 6).Value = 1806466  # F51 Volume
$ws.Cells.Item(51, 7).Value = 0.14603  # G51 Change (24h)
